# CIBMTR reporting IG - ValueSet/auto-differential-blood-vs
# Refresh the "Metadata" sheet: version bump, status -> draft, new date,
# real contact details (publisher + author), a new Jurisdiction row, a
# filled-in Description, and the Immutable row pushed down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Grow the used range by one row first: clone row 15's formatting onto the
# new row 16 so the appended row matches the existing bordered/wrapped style
# instead of picking up a blank default style.
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null

# Fill bottom-up so every row's own read happens before it gets overwritten.
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "Auto Differential panel - Blood (57023-4)"

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

$ws.Range("A10").Value = "Contact"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"
$ws.Range("B6").Value = "draft"
$ws.Range("B3").Value = "0.1.7"
